$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Purchase 22-23")

# Row 2 held a one-off opening-balance line with no invoice/date; it was
# removed and everything below it shifts up by one row.
$ws1.Rows(2).Delete()

# The shift leaves a dangling #REF! term in the subtotal formula (old row 5,
# now row 4) where the deleted row's cell used to be referenced - restore it
# to a clean SUM of the remaining three rows in the group.
$ws1.Range("F4").Formula = "=E2+E3+E4"

# The "Sr. No" sequence marker that used to sit on the deleted row's line
# (A2) needs to stay on the new first line of the group.
$ws1.Range("A2").Value = 1

# Selection now only spans the single active cell instead of the old F5:F20 block.
$ws1.Range("F5").Select() | Out-Null

# The user switched focus to the "Purchase 22-23" tab before saving.
$ws1.Activate() | Out-Null
